$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'287.54"
$ws.Range("E2").Value = "'1.61%"
$ws.Range("E3").Value = "'4.20%"
$ws.Range("D4").Value = "'5.079"
$ws.Range("E4").Value = "'1.19%"
$ws.Range("D5").Value = "'0.06763"
$ws.Range("E5").Value = "'3.97%"
$ws.Range("D6").Value = "'7.357"
$ws.Range("E6").Value = "'1.94%"
$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D7").Value = "'1.410"
$ws.Range("E7").Value = "'2.14%"
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = "'0.9145"
$ws.Range("E8").Value = "'-0.36%"
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").Value = "'0.1599"
$ws.Range("E9").Value = "'4.12%"
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = "'0.06923"
$ws.Range("E10").Value = "'8.64%"
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = "'0.07688"
$ws.Range("E11").Value = "'1.59%"
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = "'0.02919"
$ws.Range("E12").Value = "'2.05%"
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = "'0.08987"
$ws.Range("E13").Value = "'0.07%"
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = "'0.001595"
$ws.Range("E14").Value = "'0.02%"
$ws.Range("B15").Value = 'CoinExToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D15").Value = "'0.04473"
$ws.Range("E15").Value = "'0.94%"
$ws.Range("B16").Value = 'One'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D16").Value = "'0.0006464"
$ws.Range("E16").Value = "'1.43%"
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").Value = "'0.006157"
$ws.Range("E17").Value = "'-0.49%"
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").Value = "'3.451"
$ws.Range("E18").Value = "'0.17%"
$ws.Range("B19").Value = 'GateToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D19").Value = "'3.441"
$ws.Range("E19").Value = "'2.39%"
$ws.Range("D20").Value = "'2.230"
$ws.Range("E20").Value = "'-0.54%"
$ws.Range("E22").Value = "'-1.78%"
$ws.Range("D23").Value = "'4.090"
$ws.Range("E23").Value = "'3.08%"
$ws.Range("D24").Value = "'0.1580"
$ws.Range("E24").Value = "'2.38%"
$ws.Range("D25").Value = "'0.001195"
$ws.Range("E25").Value = "'1.17%"
$ws.Range("D26").Value = "'0.004145"
$ws.Range("E27").Value = "'-0.09%"
$ws.Range("E28").Value = "'-0.17%"
$ws.Range("D40").Value = "'0.04264"
$ws.Range("E40").Value = "'3.38%"
$ws.Range("D41").Value = "'0.006818"
$ws.Range("E41").Value = "'1.97%"
$ws.Range("D42").Value = "'0.1244"
$ws.Range("E42").Value = "'1.30%"
$ws.Range("D43").Value = "'0.002216"
$ws.Range("E43").Value = "'3.18%"
$ws.Range("D44").Value = "'0.01303"
$ws.Range("E44").Value = "'12.89%"
$ws.Range("D45").Value = "'0.00005685"
$ws.Range("E45").Value = "'0.58%"
$ws.Range("D46").Value = "'1.963"
$ws.Range("E46").Value = "'-0.01%"
$ws.Range("E47").Value = "'-18.74%"
